# The weekly refresh re-sorts the daily Tuna price rows (2-22).
# Read every source rows moving columns (Fecha, Volumen, Precio
# minimo/maximo/promedio, Unidad, Origen, Precio $/Kg, Kg/unidad)
# first, then write them to their new row, so rows that swap with
# each other do not clobber data still waiting to be read.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = @{}
$row2["D"] = $ws.Range("D15").Value2
$row2["M"] = $ws.Range("M15").Value2
$row2["N"] = $ws.Range("N15").Value2
$row2["O"] = $ws.Range("O15").Value2
$row2["P"] = $ws.Range("P15").Value2
$row2["Q"] = $ws.Range("Q15").Value2
$row2["R"] = $ws.Range("R15").Value2
$row2["S"] = $ws.Range("S15").Value2
$row2["T"] = $ws.Range("T15").Value2
$row3 = @{}
$row3["D"] = $ws.Range("D18").Value2
$row3["M"] = $ws.Range("M18").Value2
$row3["N"] = $ws.Range("N18").Value2
$row3["O"] = $ws.Range("O18").Value2
$row3["P"] = $ws.Range("P18").Value2
$row3["Q"] = $ws.Range("Q18").Value2
$row3["R"] = $ws.Range("R18").Value2
$row3["S"] = $ws.Range("S18").Value2
$row3["T"] = $ws.Range("T18").Value2
$row4 = @{}
$row4["D"] = $ws.Range("D12").Value2
$row4["M"] = $ws.Range("M12").Value2
$row4["N"] = $ws.Range("N12").Value2
$row4["O"] = $ws.Range("O12").Value2
$row4["P"] = $ws.Range("P12").Value2
$row4["Q"] = $ws.Range("Q12").Value2
$row4["R"] = $ws.Range("R12").Value2
$row4["S"] = $ws.Range("S12").Value2
$row4["T"] = $ws.Range("T12").Value2
$row5 = @{}
$row5["D"] = $ws.Range("D20").Value2
$row5["M"] = $ws.Range("M20").Value2
$row5["N"] = $ws.Range("N20").Value2
$row5["O"] = $ws.Range("O20").Value2
$row5["P"] = $ws.Range("P20").Value2
$row5["Q"] = $ws.Range("Q20").Value2
$row5["R"] = $ws.Range("R20").Value2
$row5["S"] = $ws.Range("S20").Value2
$row5["T"] = $ws.Range("T20").Value2
$row6 = @{}
$row6["D"] = $ws.Range("D13").Value2
$row6["M"] = $ws.Range("M13").Value2
$row6["N"] = $ws.Range("N13").Value2
$row6["O"] = $ws.Range("O13").Value2
$row6["P"] = $ws.Range("P13").Value2
$row6["Q"] = $ws.Range("Q13").Value2
$row6["R"] = $ws.Range("R13").Value2
$row6["S"] = $ws.Range("S13").Value2
$row6["T"] = $ws.Range("T13").Value2
$row7 = @{}
$row7["D"] = $ws.Range("D8").Value2
$row7["M"] = $ws.Range("M8").Value2
$row7["N"] = $ws.Range("N8").Value2
$row7["O"] = $ws.Range("O8").Value2
$row7["P"] = $ws.Range("P8").Value2
$row7["Q"] = $ws.Range("Q8").Value2
$row7["R"] = $ws.Range("R8").Value2
$row7["S"] = $ws.Range("S8").Value2
$row7["T"] = $ws.Range("T8").Value2
$row8 = @{}
$row8["D"] = $ws.Range("D14").Value2
$row8["M"] = $ws.Range("M14").Value2
$row8["N"] = $ws.Range("N14").Value2
$row8["O"] = $ws.Range("O14").Value2
$row8["P"] = $ws.Range("P14").Value2
$row8["Q"] = $ws.Range("Q14").Value2
$row8["R"] = $ws.Range("R14").Value2
$row8["S"] = $ws.Range("S14").Value2
$row8["T"] = $ws.Range("T14").Value2
$row9 = @{}
$row9["D"] = $ws.Range("D21").Value2
$row9["M"] = $ws.Range("M21").Value2
$row9["N"] = $ws.Range("N21").Value2
$row9["O"] = $ws.Range("O21").Value2
$row9["P"] = $ws.Range("P21").Value2
$row9["Q"] = $ws.Range("Q21").Value2
$row9["R"] = $ws.Range("R21").Value2
$row9["S"] = $ws.Range("S21").Value2
$row9["T"] = $ws.Range("T21").Value2
$row10 = @{}
$row10["D"] = $ws.Range("D4").Value2
$row10["M"] = $ws.Range("M4").Value2
$row10["N"] = $ws.Range("N4").Value2
$row10["O"] = $ws.Range("O4").Value2
$row10["P"] = $ws.Range("P4").Value2
$row10["Q"] = $ws.Range("Q4").Value2
$row10["R"] = $ws.Range("R4").Value2
$row10["S"] = $ws.Range("S4").Value2
$row10["T"] = $ws.Range("T4").Value2
$row11 = @{}
$row11["D"] = $ws.Range("D19").Value2
$row11["M"] = $ws.Range("M19").Value2
$row11["N"] = $ws.Range("N19").Value2
$row11["O"] = $ws.Range("O19").Value2
$row11["P"] = $ws.Range("P19").Value2
$row11["Q"] = $ws.Range("Q19").Value2
$row11["R"] = $ws.Range("R19").Value2
$row11["S"] = $ws.Range("S19").Value2
$row11["T"] = $ws.Range("T19").Value2
$row12 = @{}
$row12["D"] = $ws.Range("D6").Value2
$row12["M"] = $ws.Range("M6").Value2
$row12["N"] = $ws.Range("N6").Value2
$row12["O"] = $ws.Range("O6").Value2
$row12["P"] = $ws.Range("P6").Value2
$row12["Q"] = $ws.Range("Q6").Value2
$row12["R"] = $ws.Range("R6").Value2
$row12["S"] = $ws.Range("S6").Value2
$row12["T"] = $ws.Range("T6").Value2
$row13 = @{}
$row13["D"] = $ws.Range("D11").Value2
$row13["M"] = $ws.Range("M11").Value2
$row13["N"] = $ws.Range("N11").Value2
$row13["O"] = $ws.Range("O11").Value2
$row13["P"] = $ws.Range("P11").Value2
$row13["Q"] = $ws.Range("Q11").Value2
$row13["R"] = $ws.Range("R11").Value2
$row13["S"] = $ws.Range("S11").Value2
$row13["T"] = $ws.Range("T11").Value2
$row14 = @{}
$row14["D"] = $ws.Range("D2").Value2
$row14["M"] = $ws.Range("M2").Value2
$row14["N"] = $ws.Range("N2").Value2
$row14["O"] = $ws.Range("O2").Value2
$row14["P"] = $ws.Range("P2").Value2
$row14["Q"] = $ws.Range("Q2").Value2
$row14["R"] = $ws.Range("R2").Value2
$row14["S"] = $ws.Range("S2").Value2
$row14["T"] = $ws.Range("T2").Value2
$row15 = @{}
$row15["D"] = $ws.Range("D7").Value2
$row15["M"] = $ws.Range("M7").Value2
$row15["N"] = $ws.Range("N7").Value2
$row15["O"] = $ws.Range("O7").Value2
$row15["P"] = $ws.Range("P7").Value2
$row15["Q"] = $ws.Range("Q7").Value2
$row15["R"] = $ws.Range("R7").Value2
$row15["S"] = $ws.Range("S7").Value2
$row15["T"] = $ws.Range("T7").Value2
$row16 = @{}
$row16["D"] = $ws.Range("D10").Value2
$row16["M"] = $ws.Range("M10").Value2
$row16["N"] = $ws.Range("N10").Value2
$row16["O"] = $ws.Range("O10").Value2
$row16["P"] = $ws.Range("P10").Value2
$row16["Q"] = $ws.Range("Q10").Value2
$row16["R"] = $ws.Range("R10").Value2
$row16["S"] = $ws.Range("S10").Value2
$row16["T"] = $ws.Range("T10").Value2
$row17 = @{}
$row17["D"] = $ws.Range("D9").Value2
$row17["M"] = $ws.Range("M9").Value2
$row17["N"] = $ws.Range("N9").Value2
$row17["O"] = $ws.Range("O9").Value2
$row17["P"] = $ws.Range("P9").Value2
$row17["Q"] = $ws.Range("Q9").Value2
$row17["R"] = $ws.Range("R9").Value2
$row17["S"] = $ws.Range("S9").Value2
$row17["T"] = $ws.Range("T9").Value2
$row18 = @{}
$row18["D"] = $ws.Range("D3").Value2
$row18["M"] = $ws.Range("M3").Value2
$row18["N"] = $ws.Range("N3").Value2
$row18["O"] = $ws.Range("O3").Value2
$row18["P"] = $ws.Range("P3").Value2
$row18["Q"] = $ws.Range("Q3").Value2
$row18["R"] = $ws.Range("R3").Value2
$row18["S"] = $ws.Range("S3").Value2
$row18["T"] = $ws.Range("T3").Value2
$row19 = @{}
$row19["D"] = $ws.Range("D5").Value2
$row19["M"] = $ws.Range("M5").Value2
$row19["N"] = $ws.Range("N5").Value2
$row19["O"] = $ws.Range("O5").Value2
$row19["P"] = $ws.Range("P5").Value2
$row19["Q"] = $ws.Range("Q5").Value2
$row19["R"] = $ws.Range("R5").Value2
$row19["S"] = $ws.Range("S5").Value2
$row19["T"] = $ws.Range("T5").Value2
$row20 = @{}
$row20["D"] = $ws.Range("D22").Value2
$row20["M"] = $ws.Range("M22").Value2
$row20["N"] = $ws.Range("N22").Value2
$row20["O"] = $ws.Range("O22").Value2
$row20["P"] = $ws.Range("P22").Value2
$row20["Q"] = $ws.Range("Q22").Value2
$row20["R"] = $ws.Range("R22").Value2
$row20["S"] = $ws.Range("S22").Value2
$row20["T"] = $ws.Range("T22").Value2
$row21 = @{}
$row21["D"] = $ws.Range("D16").Value2
$row21["M"] = $ws.Range("M16").Value2
$row21["N"] = $ws.Range("N16").Value2
$row21["O"] = $ws.Range("O16").Value2
$row21["P"] = $ws.Range("P16").Value2
$row21["Q"] = $ws.Range("Q16").Value2
$row21["R"] = $ws.Range("R16").Value2
$row21["S"] = $ws.Range("S16").Value2
$row21["T"] = $ws.Range("T16").Value2
$row22 = @{}
$row22["D"] = $ws.Range("D17").Value2
$row22["M"] = $ws.Range("M17").Value2
$row22["N"] = $ws.Range("N17").Value2
$row22["O"] = $ws.Range("O17").Value2
$row22["P"] = $ws.Range("P17").Value2
$row22["Q"] = $ws.Range("Q17").Value2
$row22["R"] = $ws.Range("R17").Value2
$row22["S"] = $ws.Range("S17").Value2
$row22["T"] = $ws.Range("T17").Value2

# Now write the captured values into their destination rows.
$ws.Range("D2").Value = $row2["D"]
$ws.Range("M2").Value = $row2["M"]
$ws.Range("N2").Value = $row2["N"]
$ws.Range("O2").Value = $row2["O"]
$ws.Range("P2").Value = $row2["P"]
$ws.Range("Q2").Value = $row2["Q"]
$ws.Range("R2").Value = $row2["R"]
$ws.Range("S2").Value = $row2["S"]
$ws.Range("T2").Value = $row2["T"]
$ws.Range("D3").Value = $row3["D"]
$ws.Range("M3").Value = $row3["M"]
$ws.Range("N3").Value = $row3["N"]
$ws.Range("O3").Value = $row3["O"]
$ws.Range("P3").Value = $row3["P"]
$ws.Range("Q3").Value = $row3["Q"]
$ws.Range("R3").Value = $row3["R"]
$ws.Range("S3").Value = $row3["S"]
$ws.Range("T3").Value = $row3["T"]
$ws.Range("D4").Value = $row4["D"]
$ws.Range("M4").Value = $row4["M"]
$ws.Range("N4").Value = $row4["N"]
$ws.Range("O4").Value = $row4["O"]
$ws.Range("P4").Value = $row4["P"]
$ws.Range("Q4").Value = $row4["Q"]
$ws.Range("R4").Value = $row4["R"]
$ws.Range("S4").Value = $row4["S"]
$ws.Range("T4").Value = $row4["T"]
$ws.Range("D5").Value = $row5["D"]
$ws.Range("M5").Value = $row5["M"]
$ws.Range("N5").Value = $row5["N"]
$ws.Range("O5").Value = $row5["O"]
$ws.Range("P5").Value = $row5["P"]
$ws.Range("Q5").Value = $row5["Q"]
$ws.Range("R5").Value = $row5["R"]
$ws.Range("S5").Value = $row5["S"]
$ws.Range("T5").Value = $row5["T"]
$ws.Range("D6").Value = $row6["D"]
$ws.Range("M6").Value = $row6["M"]
$ws.Range("N6").Value = $row6["N"]
$ws.Range("O6").Value = $row6["O"]
$ws.Range("P6").Value = $row6["P"]
$ws.Range("Q6").Value = $row6["Q"]
$ws.Range("R6").Value = $row6["R"]
$ws.Range("S6").Value = $row6["S"]
$ws.Range("T6").Value = $row6["T"]
$ws.Range("D7").Value = $row7["D"]
$ws.Range("M7").Value = $row7["M"]
$ws.Range("N7").Value = $row7["N"]
$ws.Range("O7").Value = $row7["O"]
$ws.Range("P7").Value = $row7["P"]
$ws.Range("Q7").Value = $row7["Q"]
$ws.Range("R7").Value = $row7["R"]
$ws.Range("S7").Value = $row7["S"]
$ws.Range("T7").Value = $row7["T"]
$ws.Range("D8").Value = $row8["D"]
$ws.Range("M8").Value = $row8["M"]
$ws.Range("N8").Value = $row8["N"]
$ws.Range("O8").Value = $row8["O"]
$ws.Range("P8").Value = $row8["P"]
$ws.Range("Q8").Value = $row8["Q"]
$ws.Range("R8").Value = $row8["R"]
$ws.Range("S8").Value = $row8["S"]
$ws.Range("T8").Value = $row8["T"]
$ws.Range("D9").Value = $row9["D"]
$ws.Range("M9").Value = $row9["M"]
$ws.Range("N9").Value = $row9["N"]
$ws.Range("O9").Value = $row9["O"]
$ws.Range("P9").Value = $row9["P"]
$ws.Range("Q9").Value = $row9["Q"]
$ws.Range("R9").Value = $row9["R"]
$ws.Range("S9").Value = $row9["S"]
$ws.Range("T9").Value = $row9["T"]
$ws.Range("D10").Value = $row10["D"]
$ws.Range("M10").Value = $row10["M"]
$ws.Range("N10").Value = $row10["N"]
$ws.Range("O10").Value = $row10["O"]
$ws.Range("P10").Value = $row10["P"]
$ws.Range("Q10").Value = $row10["Q"]
$ws.Range("R10").Value = $row10["R"]
$ws.Range("S10").Value = $row10["S"]
$ws.Range("T10").Value = $row10["T"]
$ws.Range("D11").Value = $row11["D"]
$ws.Range("M11").Value = $row11["M"]
$ws.Range("N11").Value = $row11["N"]
$ws.Range("O11").Value = $row11["O"]
$ws.Range("P11").Value = $row11["P"]
$ws.Range("Q11").Value = $row11["Q"]
$ws.Range("R11").Value = $row11["R"]
$ws.Range("S11").Value = $row11["S"]
$ws.Range("T11").Value = $row11["T"]
$ws.Range("D12").Value = $row12["D"]
$ws.Range("M12").Value = $row12["M"]
$ws.Range("N12").Value = $row12["N"]
$ws.Range("O12").Value = $row12["O"]
$ws.Range("P12").Value = $row12["P"]
$ws.Range("Q12").Value = $row12["Q"]
$ws.Range("R12").Value = $row12["R"]
$ws.Range("S12").Value = $row12["S"]
$ws.Range("T12").Value = $row12["T"]
$ws.Range("D13").Value = $row13["D"]
$ws.Range("M13").Value = $row13["M"]
$ws.Range("N13").Value = $row13["N"]
$ws.Range("O13").Value = $row13["O"]
$ws.Range("P13").Value = $row13["P"]
$ws.Range("Q13").Value = $row13["Q"]
$ws.Range("R13").Value = $row13["R"]
$ws.Range("S13").Value = $row13["S"]
$ws.Range("T13").Value = $row13["T"]
$ws.Range("D14").Value = $row14["D"]
$ws.Range("M14").Value = $row14["M"]
$ws.Range("N14").Value = $row14["N"]
$ws.Range("O14").Value = $row14["O"]
$ws.Range("P14").Value = $row14["P"]
$ws.Range("Q14").Value = $row14["Q"]
$ws.Range("R14").Value = $row14["R"]
$ws.Range("S14").Value = $row14["S"]
$ws.Range("T14").Value = $row14["T"]
$ws.Range("D15").Value = $row15["D"]
$ws.Range("M15").Value = $row15["M"]
$ws.Range("N15").Value = $row15["N"]
$ws.Range("O15").Value = $row15["O"]
$ws.Range("P15").Value = $row15["P"]
$ws.Range("Q15").Value = $row15["Q"]
$ws.Range("R15").Value = $row15["R"]
$ws.Range("S15").Value = $row15["S"]
$ws.Range("T15").Value = $row15["T"]
$ws.Range("D16").Value = $row16["D"]
$ws.Range("M16").Value = $row16["M"]
$ws.Range("N16").Value = $row16["N"]
$ws.Range("O16").Value = $row16["O"]
$ws.Range("P16").Value = $row16["P"]
$ws.Range("Q16").Value = $row16["Q"]
$ws.Range("R16").Value = $row16["R"]
$ws.Range("S16").Value = $row16["S"]
$ws.Range("T16").Value = $row16["T"]
$ws.Range("D17").Value = $row17["D"]
$ws.Range("M17").Value = $row17["M"]
$ws.Range("N17").Value = $row17["N"]
$ws.Range("O17").Value = $row17["O"]
$ws.Range("P17").Value = $row17["P"]
$ws.Range("Q17").Value = $row17["Q"]
$ws.Range("R17").Value = $row17["R"]
$ws.Range("S17").Value = $row17["S"]
$ws.Range("T17").Value = $row17["T"]
$ws.Range("D18").Value = $row18["D"]
$ws.Range("M18").Value = $row18["M"]
$ws.Range("N18").Value = $row18["N"]
$ws.Range("O18").Value = $row18["O"]
$ws.Range("P18").Value = $row18["P"]
$ws.Range("Q18").Value = $row18["Q"]
$ws.Range("R18").Value = $row18["R"]
$ws.Range("S18").Value = $row18["S"]
$ws.Range("T18").Value = $row18["T"]
$ws.Range("D19").Value = $row19["D"]
$ws.Range("M19").Value = $row19["M"]
$ws.Range("N19").Value = $row19["N"]
$ws.Range("O19").Value = $row19["O"]
$ws.Range("P19").Value = $row19["P"]
$ws.Range("Q19").Value = $row19["Q"]
$ws.Range("R19").Value = $row19["R"]
$ws.Range("S19").Value = $row19["S"]
$ws.Range("T19").Value = $row19["T"]
$ws.Range("D20").Value = $row20["D"]
$ws.Range("M20").Value = $row20["M"]
$ws.Range("N20").Value = $row20["N"]
$ws.Range("O20").Value = $row20["O"]
$ws.Range("P20").Value = $row20["P"]
$ws.Range("Q20").Value = $row20["Q"]
$ws.Range("R20").Value = $row20["R"]
$ws.Range("S20").Value = $row20["S"]
$ws.Range("T20").Value = $row20["T"]
$ws.Range("D21").Value = $row21["D"]
$ws.Range("M21").Value = $row21["M"]
$ws.Range("N21").Value = $row21["N"]
$ws.Range("O21").Value = $row21["O"]
$ws.Range("P21").Value = $row21["P"]
$ws.Range("Q21").Value = $row21["Q"]
$ws.Range("R21").Value = $row21["R"]
$ws.Range("S21").Value = $row21["S"]
$ws.Range("T21").Value = $row21["T"]
$ws.Range("D22").Value = $row22["D"]
$ws.Range("M22").Value = $row22["M"]
$ws.Range("N22").Value = $row22["N"]
$ws.Range("O22").Value = $row22["O"]
$ws.Range("P22").Value = $row22["P"]
$ws.Range("Q22").Value = $row22["Q"]
$ws.Range("R22").Value = $row22["R"]
$ws.Range("S22").Value = $row22["S"]
$ws.Range("T22").Value = $row22["T"]
